$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: new tag entry "天生丽质" (Beautiful by nature)
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = "天生丽质"
$ws.Range("D4").Value = "好漂亮的姐姐~"
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = $true
$ws.Range("G4").Value = $true
$ws.Range("H4").Value = $true
$ws.Range("I4").Value = $true

# Update the active selection to I4
$ws.Range("I4").Select()
